$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-24 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-25 Thursday", 2) | Out-Null
$d.Content.Find.Execute("70×68=", $true, $false, $false, $false, $false, $true, 1, $false, "71×69=", 2) | Out-Null
$d.Content.Find.Execute("21×43=", $true, $false, $false, $false, $false, $true, 1, $false, "31×77=", 2) | Out-Null
$d.Content.Find.Execute("36×43=", $true, $false, $false, $false, $false, $true, 1, $false, "37×59=", 2) | Out-Null
$d.Content.Find.Execute("48×61=", $true, $false, $false, $false, $false, $true, 1, $false, "86×99=", 2) | Out-Null
$d.Content.Find.Execute("25×46=", $true, $false, $false, $false, $false, $true, 1, $false, "37×92=", 2) | Out-Null
$d.Content.Find.Execute("88×47=", $true, $false, $false, $false, $false, $true, 1, $false, "38×55=", 2) | Out-Null
$d.Content.Find.Execute("32×84=", $true, $false, $false, $false, $false, $true, 1, $false, "76×38=", 2) | Out-Null
$d.Content.Find.Execute("96×74=", $true, $false, $false, $false, $false, $true, 1, $false, "84×99=", 2) | Out-Null
$d.Content.Find.Execute("27×49=", $true, $false, $false, $false, $false, $true, 1, $false, "84×40=", 2) | Out-Null
$d.Content.Find.Execute("48×40=", $true, $false, $false, $false, $false, $true, 1, $false, "55×55=", 2) | Out-Null
$d.Content.Find.Execute("93×88=", $true, $false, $false, $false, $false, $true, 1, $false, "29×12=", 2) | Out-Null
$d.Content.Find.Execute("68×60=", $true, $false, $false, $false, $false, $true, 1, $false, "87×80=", 2) | Out-Null
$d.Content.Find.Execute("87×62=", $true, $false, $false, $false, $false, $true, 1, $false, "37×41=", 2) | Out-Null
$d.Content.Find.Execute("22×43=", $true, $false, $false, $false, $false, $true, 1, $false, "78×13=", 2) | Out-Null
$d.Content.Find.Execute("74×59=", $true, $false, $false, $false, $false, $true, 1, $false, "54×64=", 2) | Out-Null
$d.Content.Find.Execute("82×52=", $true, $false, $false, $false, $false, $true, 1, $false, "81×79=", 2) | Out-Null
$d.Content.Find.Execute("88×35=", $true, $false, $false, $false, $false, $true, 1, $false, "85×66=", 2) | Out-Null
$d.Content.Find.Execute("94×80=", $true, $false, $false, $false, $false, $true, 1, $false, "57×55=", 2) | Out-Null
$d.Content.Find.Execute("39×25=", $true, $false, $false, $false, $false, $true, 1, $false, "76×71=", 2) | Out-Null
$d.Content.Find.Execute("26×95=", $true, $false, $false, $false, $false, $true, 1, $false, "57×51=", 2) | Out-Null
$d.Content.Find.Execute("46×52=", $true, $false, $false, $false, $false, $true, 1, $false, "31×77=", 2) | Out-Null
$d.Content.Find.Execute("66×16=", $true, $false, $false, $false, $false, $true, 1, $false, "61×26=", 2) | Out-Null
$d.Content.Find.Execute("11×14=", $true, $false, $false, $false, $false, $true, 1, $false, "43×72=", 2) | Out-Null
$d.Content.Find.Execute("35×46=", $true, $false, $false, $false, $false, $true, 1, $false, "79×51=", 2) | Out-Null
$d.Content.Find.Execute("79×66=", $true, $false, $false, $false, $false, $true, 1, $false, "88×61=", 2) | Out-Null
